$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 158; existing rows 158:283 shift down to 159:284.
$ws.Rows("158:158").Insert()

# Populate the newly inserted row 158 with the new weekly data point.
$ws.Cells.Item(158, 1).Value = 7
$ws.Cells.Item(158, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(158, 3).Value = "Ñuble"
$ws.Cells.Item(158, 4).Value = 44673
$ws.Cells.Item(158, 5).Value = 16
$ws.Cells.Item(158, 6).Value = 100114013
$ws.Cells.Item(158, 7).Value = "Zanahoria"
$ws.Cells.Item(158, 8).Value = "Sin especificar"
$ws.Cells.Item(158, 9).Value = "Primera"
$ws.Cells.Item(158, 10).Value = 120
$ws.Cells.Item(158, 11).Value = 6500
$ws.Cells.Item(158, 12).Value = 7000
$ws.Cells.Item(158, 13).Value = 6750
$ws.Cells.Item(158, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(158, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(158, 16).Value = 338
$ws.Cells.Item(158, 17).Value = 20
$ws.Cells.Item(158, 18).Value = "Hortaliza"
